$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Rows.Item(11).Insert()
Write-Host "inserted"
for ($i = 1; $i -le 20; $i++) {
    $a = $ws.Cells.Item($i, 1).Value()
    $b = $ws.Cells.Item($i, 2).Value()
    Write-Host "$i : A=$a | B=$b"
}
